# Swap the presentation's theme palette from the custom "Integral" /
# "Red Violet" scheme back to the stock PowerPoint "Office Theme" /
# "Office" color scheme (theme1.xml).
#
# Office theme RGB values (as hex RRGGBB), converted to the BGR-packed
# long that PowerPoint's ThemeColor.RGB property expects:
#   dk1      000000 -> 0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A -> 6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 -> 3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 -> 49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 -> 4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 -> 7491477

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

$officeColors = @(
    0,
    16777215,
    6968388,
    15132391,
    13998939,
    3243501,
    10855845,
    49407,
    12874308,
    4697456,
    12673797,
    7491477
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
